$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price row (2023-04-05, serial 45021) was added to the
# "Espinaca" price series for Vega Central Mapocho de Santiago. It is
# inserted immediately above the row that used to be row 546 (which,
# together with every row below it through the former row 586, shifts
# down by one row).
$ws.Rows(546).Insert()

$ws.Cells.Item(546, 1).Value = 9
$ws.Cells.Item(546, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(546, 3).Value = "Metropolitana"
$ws.Cells.Item(546, 4).Value = 45021
$ws.Cells.Item(546, 5).Value = 13
$ws.Cells.Item(546, 6).Value = 100112012
$ws.Cells.Item(546, 7).Value = "Espinaca"
$ws.Cells.Item(546, 8).Value = "Sin especificar"
$ws.Cells.Item(546, 9).Value = "Primera"
$ws.Cells.Item(546, 10).Value = 240
$ws.Cells.Item(546, 11).Value = 9000
$ws.Cells.Item(546, 12).Value = 10000
$ws.Cells.Item(546, 13).Value = 9417
$ws.Cells.Item(546, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(546, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(546, 16).Value = 942
$ws.Cells.Item(546, 17).Value = 10
$ws.Cells.Item(546, 18).Value = "Hortaliza"
